# Sync attendance_reports content: correct the "Recorded By" column (G)
# so that the ordering of "System" / "system" and the recorder's e-mail
# address matches the canonical data from the main repository.
#
# Rule observed in the source data:
#   "System, dnasr281@gmail.com"            -> "dnasr281@gmail.com, System"
#   "System, system, backup@backdoor.com"   -> "System, backup@backdoor.com, system"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count() + $usedRange.Row() - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value()

    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
    elseif ($val -eq "System, system, backup@backdoor.com") {
        $cell.Value = "System, backup@backdoor.com, system"
    }
}
